$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking strings (e.g. with dots as
# thousands/decimal separators) are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.444.77'
$ws.Range("E2").Value = '  -0.89%  '

$ws.Range("D3").Value = '1.823.56'
$ws.Range("E3").Value = '  -2.22%  '

$ws.Range("E4").Value = '  -0.67%  '

$ws.Range("D5").Value = '332.38'
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("E6").Value = '  -0.55%  '

$ws.Range("D7").Value = '0.4569'
$ws.Range("E7").Value = '  -2.85%  '

$ws.Range("D8").Value = '0.3794'
$ws.Range("E8").Value = '  -3.32%  '

$ws.Range("D9").Value = '46.22'
$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("D10").Value = '0.07879'
$ws.Range("E10").Value = '  -1.54%  '

$ws.Range("D11").Value = '0.9676'
$ws.Range("E11").Value = '  -3.75%  '

$ws.Range("D12").Value = '20.98'
$ws.Range("E12").Value = '  -3.84%  '

$ws.Range("D13").Value = '5.873'
$ws.Range("E13").Value = '  -2.08%  '

$ws.Range("D14").Value = '1.814.83'
$ws.Range("E14").Value = '  -3.78%  '

$ws.Range("D15").Value = '7.035'
$ws.Range("E15").Value = '  -3.07%  '

$ws.Range("E16").Value = '  -0.58%  '

$ws.Range("D17").Value = '89.59'
$ws.Range("E17").Value = '  +1.35%  '

$ws.Range("D18").Value = '0.06640'
$ws.Range("E18").Value = '  -1.33%  '

$ws.Range("D19").Value = '0.00001024'
$ws.Range("E19").Value = '  -1.94%  '

$ws.Range("D20").Value = '17.06'
$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("E21").Value = '  -0.62%  '

$ws.Range("D22").Value = '27.420.92'
$ws.Range("E22").Value = '  -0.86%  '

$ws.Range("D23").Value = '5.323'
$ws.Range("E23").Value = '  -2.68%  '

$ws.Range("D24").Value = '10.79'
$ws.Range("E24").Value = '  -1.13%  '

$ws.Range("D25").Value = '2.308'
$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").Value = '2.031.44'
$ws.Range("E26").Value = '  -3.38%  '

$ws.Range("D27").Value = '155.47'
$ws.Range("E27").Value = '  -2.53%  '

$ws.Range("D28").Value = '19.34'
$ws.Range("E28").Value = '  -2.52%  '

$ws.Range("D29").Value = '2.048'
$ws.Range("E29").Value = '  -4.52%  '

$ws.Range("D30").Value = '5.269'
$ws.Range("E30").Value = '  -3.23%  '

$ws.Range("D31").Value = '118.07'
$ws.Range("E31").Value = '  -3.01%  '

$ws.Range("D32").Value = '0.9391'
$ws.Range("E32").Value = '  -4.39%  '

$ws.Range("D33").Value = '0.09291'
$ws.Range("E33").Value = '  -2.25%  '

$ws.Range("D34").Value = '3.580'

$ws.Range("D35").Value = '5.237'
$ws.Range("E35").Value = '  -1.24%  '

$ws.Range("D36").Value = '1.315'
$ws.Range("E36").Value = '  -1.63%  '

$ws.Range("D37").Value = '0.05923'
$ws.Range("E37").Value = '  -2.38%  '

$ws.Range("D38").Value = '0.02177'
$ws.Range("E38").Value = '  -2.44%  '

$ws.Range("D39").Value = '8.036'
$ws.Range("E39").Value = '  -3.45%  '

$ws.Range("D40").Value = '1.146'
$ws.Range("E40").Value = '  -4.17%  '

$ws.Range("D41").Value = '0.5752'
$ws.Range("E41").Value = '  -3.58%  '

$ws.Range("D42").Value = '0.1821'
$ws.Range("E42").Value = '  -3.47%  '

$ws.Range("D43").Value = '9.938'
$ws.Range("E43").Value = '  -3.39%  '

$ws.Range("D44").Value = '1.283'
$ws.Range("E44").Value = '  +2.51%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '11.97'
$ws.Range("E45").Value = '  -1.78%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.5435'
$ws.Range("E46").Value = '  -3.75%  '

$ws.Range("D47").Value = '1.861'
$ws.Range("E47").Value = '  -3.44%  '

$ws.Range("D48").Value = '110.55'
$ws.Range("E48").Value = '  -1.48%  '

$ws.Range("D49").Value = '0.06596'
$ws.Range("E49").Value = '  -2.56%  '

$ws.Range("D50").Value = '1.004'
$ws.Range("E50").Value = '  -0.69%  '

$ws.Range("E51").Value = '  -1.68%  '
